$d = $word.ActiveDocument
$brk = [string][char]11

# --- 1. Remove the stray manual line break after
#        "cancers have not spread to the lymph nodes"
#        (surgical delete of just the <w:br/> char so the surrounding
#        runs/formatting are left untouched). ---
$rng = $d.Content
if ($rng.Find.Execute("cancers have not spread to the lymph nodes" + $brk, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $brkRange = $d.Range($rng.End - 1, $rng.End)
    $brkRange.Delete()
}

# --- 2. Remove the stray manual line break after
#        "Chemotherapy + radiation given together over 6 weeks" ---
$rng = $d.Content
if ($rng.Find.Execute("Chemotherapy + radiation given together over 6 weeks" + $brk, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $brkRange = $d.Range($rng.End - 1, $rng.End)
    $brkRange.Delete()
}

# --- 3. Remove the stand-alone "First Paragraph" paragraph (it contains only a
#        manual line break) sitting between the "34 My Atrium Patient Portal"
#        heading and the "Critical to good communication..." bullet. ---
$paras = $d.Paragraphs
$n = $paras.Count
$target = $null
for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -eq ($brk + [string][char]13)) {
        $styleName = $p.Style.NameLocal
        if ($styleName -eq "First Paragraph") {
            $prev = $p.Previous()
            if ($prev.Range.Text -like "*My Atrium Patient Portal*") {
                $target = $p
            }
        }
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# --- 4. Wording tweaks ---
$d.Content.Find.Execute("Critical to good communication with your cancer care team", $true, $false, $false, $false, $false, $true, 1, $false, "Critical to good communication with your care team", 2) | Out-Null

$d.Content.Find.Execute("Important to reduce the risk of complications from cancer treatment", $true, $false, $false, $false, $false, $true, 1, $false, "Reduces risk of complications from treatment", 2) | Out-Null

$d.Content.Find.Execute("Working hard enough that you can’t carry a conversation", $true, $false, $false, $false, $false, $true, 1, $false, "Working hard enough that you can’t converse", 2) | Out-Null

$d.Content.Find.Execute("Start slow an build up", $true, $false, $false, $false, $false, $true, 1, $false, "Start slowly and build up", 2) | Out-Null

$d.Content.Find.Execute("Smoking makes it more difficult to get through cancer treatment", $true, $false, $false, $false, $false, $true, 1, $false, "Smoking makes cancer treatment more difficult", 2) | Out-Null

$d.Content.Find.Execute("American Lung Asssociation fredomfromsmoking.org", $true, $false, $false, $false, $false, $true, 1, $false, "American Lung Assn fredomfromsmoking.org", 2) | Out-Null

$d.Content.Find.Execute("1:1 Smoking Cessation Counseling Clinics (Metro Charlotte)", $true, $false, $false, $false, $false, $true, 1, $false, "1:1 Smoking Cessation Counseling (Metro Charlotte)", 2) | Out-Null
